$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New column E: "Költség becslés 20 kérdésre" (Estimated cost of 20 questions) ---

# Header cell E3: same look as the other header cells (copy format from D3), then set text.
$ws.Range("D3").Copy()
$ws.Range("E3").PasteSpecial(-4122)
$ws.Range("E3").Value = "Költség becslés 20 kérdésre"

# Body cells that have no cost estimate get a plain "-" like D-column's normal style (copied from D4).
$dashCells = @("E4", "E6", "E7", "E8", "E13")
foreach ($addr in $dashCells) {
    $ws.Range("D4").Copy()
    $ws.Range($addr).PasteSpecial(-4122)
    $ws.Range($addr).Value = "-"
}

# Body cells that do have a numeric cost estimate: same base style as D4, but with a currency number format.
$costValues = @{ "E5" = 0.2; "E9" = 1.3; "E10" = 1.29; "E11" = 6.5; "E12" = 0.5 }
foreach ($addr in $costValues.Keys) {
    $ws.Range("D4").Copy()
    $ws.Range($addr).PasteSpecial(-4122)
    $ws.Range($addr).Value = $costValues[$addr]
    $ws.Range($addr).NumberFormat = "$#,##0.00_);[Red]($#,##0.00)"
}

$excel.CutCopyMode = 0

# Column E should be as wide as column D.
$ws.Range("E1").ColumnWidth = $ws.Range("D1").ColumnWidth

# Restore the selection to where it ended up after the edit.
$ws.Range("C18").Select()

Write-Output "done"
